$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.615.14"
$ws.Range("E2").Value = "  -2.70%  "
$ws.Range("D3").Value = "3.708.13"
$ws.Range("E3").Value = "  -3.47%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'600.62"
$ws.Range("E5").Value = "  +1.25%  "
$ws.Range("D6").Value = "'184.69"
$ws.Range("E6").Value = "  +11.38%  "
$ws.Range("D7").Value = "3.701.85"
$ws.Range("E7").Value = "  -3.47%  "
$ws.Range("D8").Value = "'0.633"
$ws.Range("E8").Value = "  -5.53%  "
$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = "  +0.15%  "
$ws.Range("D10").Value = "'0.720"
$ws.Range("D11").Value = "'0.164"
$ws.Range("D12").Value = "'56.85"
$ws.Range("E12").Value = "  +7.46%  "
$ws.Range("D14").Value = "'10.46"
$ws.Range("E14").Value = "  -7.60%  "
$ws.Range("D15").Value = "4.296.20"
$ws.Range("E15").Value = "  -3.60%  "
$ws.Range("D16").Value = "3.709.33"
$ws.Range("E16").Value = "  -3.82%  "
$ws.Range("D17").Value = "'19.41"
$ws.Range("E17").Value = "  -7.43%  "
$ws.Range("E18").Value = "  -2.06%  "
$ws.Range("E19").Value = "  -6.37%  "
$ws.Range("E20").Value = "  -6.51%  "
$ws.Range("D21").Value = "68.332.28"
$ws.Range("E21").Value = "  -2.98%  "
$ws.Range("D22").Value = "'411.55"
$ws.Range("E22").Value = "  -5.48%  "
$ws.Range("E23").Value = "  -1.17%  "
$ws.Range("D24").Value = "'89.40"
$ws.Range("E24").Value = "  -4.56%  "
$ws.Range("E25").Value = "  -6.65%  "
$ws.Range("D26").Value = "'12.87"
$ws.Range("E26").Value = "  -7.03%  "
$ws.Range("D27").Value = "'10.93"
$ws.Range("E27").Value = "  -2.12%  "
$ws.Range("D28").Value = "'3.93"
$ws.Range("E28").Value = "  -0.57%  "
$ws.Range("E29").Value = "  +1.78%  "
$ws.Range("D30").Value = "'9.51"
$ws.Range("E30").Value = "  -8.75%  "
$ws.Range("D31").Value = "'32.96"
$ws.Range("E31").Value = "  -5.85%  "
$ws.Range("D32").Value = "'7.26"
$ws.Range("E32").Value = "  -10.22%  "
$ws.Range("E33").Value = "  -6.79%  "
$ws.Range("D34").Value = "'0.117"
$ws.Range("E34").Value = "  -5.72%  "
$ws.Range("D35").Value = "'43.91"
$ws.Range("E35").Value = "  -8.63%  "
$ws.Range("D36").Value = "'64.94"
$ws.Range("E36").Value = "  -6.33%  "
$ws.Range("D37").Value = "'607.18"
$ws.Range("E37").Value = "  -4.36%  "
$ws.Range("E38").Value = "  -9.08%  "
$ws.Range("E39").Value = "  -5.14%  "
$ws.Range("E40").Value = "  +0.14%  "
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("E42").Value = "  -5.68%  "
$ws.Range("D43").Value = "'3.05"
$ws.Range("E43").Value = "  -6.37%  "
$ws.Range("E44").Value = "  +1.16%  "
$ws.Range("E45").Value = "  -5.60%  "
$ws.Range("E46").Value = "  -9.99%  "
$ws.Range("D47").Value = "'9.26"
$ws.Range("E47").Value = "  -6.88%  "
$ws.Range("D48").Value = "'2.74"
$ws.Range("E48").Value = "  -3.67%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "2.793.33"
$ws.Range("E49").Value = "  -1.82%  "
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").Value = "'0.135"
$ws.Range("E50").Value = "  -5.65%  "
$ws.Range("E51").Value = "  -2.26%  "
